# Update the "La Peñita" monitoring data: new counts were uploaded and the
# rows re-sorted in descending order by total_registros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Name = "PEREZ VEGA ANA YSABEL";          Total = 107 },
    @{ Row = 3;  Name = "ZAPATA ZETA ROSA ARACELI";       Total = 101 },
    @{ Row = 4;  Name = "TIMOTEO BAYONA SHARYN LISSETH";  Total = 99 },
    @{ Row = 5;  Name = "PANTA MONZON SHIRLEY MARIBEL";   Total = 89 },
    @{ Row = 6;  Name = "GARAVITO LEON IVONNE LISSETH";   Total = 79 },
    @{ Row = 7;  Name = "VALLE SILVA SUTMMER ORFELINDA";  Total = 71 },
    @{ Row = 8;  Name = "NIÑO GUERRERO ANYELA MELINA";    Total = 70 },
    @{ Row = 9;  Name = "TIZON NUÑEZ FRESIA YAMILI";      Total = 70 },
    @{ Row = 10; Name = "CASTRO JUAREZ MARIA ISABEL";     Total = 66 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Total
}
